$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 437941.8
$ws.Range("J17").Value = 437941.8
$ws.Range("L17").Value = 1313825.4
$ws.Range("N17").Value = -1314161.4

# Row 107
$ws.Range("H107").Value = 1051.5186
$ws.Range("I107").Value = 625.8421
$ws.Range("J107").Value = 2062.5
$ws.Range("K107").Value = 625.8421
$ws.Range("L107").Value = 2062.5
$ws.Range("M107").Value = 1294.1579
$ws.Range("N107").Value = -5902.5

# Row 116
$ws.Range("H116").Value = 3400
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

# Row 138
$ws.Range("H138").Value = 135473.62
$ws.Range("I138").Value = 22408.479
$ws.Range("J138").Value = 279945.75
$ws.Range("K138").Value = 67225.43700000001
$ws.Range("L138").Value = 839837.25
$ws.Range("M138").Value = -62085.43700000001
$ws.Range("N138").Value = -850117.25

# Row 141
$ws.Range("H141").Value = 960.8
$ws.Range("I141").Value = 762.3333
$ws.Range("J141").Value = 2250.8333
$ws.Range("K141").Value = 2286.9999
$ws.Range("L141").Value = 6752.499899999999
$ws.Range("M141").Value = 2893.0001
$ws.Range("N141").Value = -17112.4999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2814.0708
$ws.Range("I32").Value = 2719.0745
$ws.Range("K32").Value = 2719.0745
$ws.Range("M32").Value = -2432.0745

# Row 74
$ws.Range("H74").Value = 1891.6538
$ws.Range("I74").Value = 1303.0476
$ws.Range("K74").Value = 1303.0476
$ws.Range("M74").Value = -429.0476000000001

# Row 77
$ws.Range("H77").Value = 1891.6538
$ws.Range("I77").Value = 1303.0476
$ws.Range("K77").Value = 6515.238
$ws.Range("M77").Value = -2147.238

# Row 131
$ws.Range("H131").Value = 57147.9
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 57147.9
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 57147.9
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -67227.89999999999

# Row 132
$ws.Range("H132").Value = 4444.149
$ws.Range("I132").Value = 2613.95
$ws.Range("K132").Value = 7841.849999999999
$ws.Range("M132").Value = -5311.849999999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 97
$ws.Range("H97").Value = 11905.818
$ws.Range("I97").Value = 4675.5557
$ws.Range("J97").Value = 44442
$ws.Range("K97").Value = 4675.5557
$ws.Range("L97").Value = 44442
$ws.Range("M97").Value = -3684.5557
$ws.Range("N97").Value = -46424

# Row 113
$ws.Range("H113").Value = 6000
$ws.Range("I113").Value = 6000
$ws.Range("K113").Value = 6000
$ws.Range("M113").Value = -3830

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2451.4546
$ws.Range("I31").Value = 1791.2122
$ws.Range("K31").Value = 1791.2122
$ws.Range("M31").Value = -1496.2122

# Row 34
$ws.Range("H34").Value = 2451.4546
$ws.Range("I34").Value = 1791.2122
$ws.Range("K34").Value = 1791.2122
$ws.Range("M34").Value = -1589.2122

# Row 58
$ws.Range("H58").Value = 1861
$ws.Range("I58").Value = 2040.2
$ws.Range("J58").Value = 1263.6666
$ws.Range("K58").Value = 2040.2
$ws.Range("L58").Value = 1263.6666
$ws.Range("M58").Value = -1837.2
$ws.Range("N58").Value = -1669.6666

# Row 134
$ws.Range("H134").Value = 2425.2856
$ws.Range("I134").Value = 2077.303
$ws.Range("J134").Value = 3701.2222
$ws.Range("K134").Value = 6231.909
$ws.Range("L134").Value = 11103.6666
$ws.Range("M134").Value = -3696.909
$ws.Range("N134").Value = -16173.6666

# Row 136
$ws.Range("H136").Value = 1861
$ws.Range("I136").Value = 2040.2
$ws.Range("J136").Value = 1263.6666
$ws.Range("K136").Value = 6120.6
$ws.Range("L136").Value = 3790.9998
$ws.Range("M136").Value = -3570.6
$ws.Range("N136").Value = -8890.9998

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 29309.215
$ws.Range("I6").Value = 34110
$ws.Range("K6").Value = 102330
$ws.Range("M6").Value = -102217

# Row 98
$ws.Range("H98").Value = 910.6923
$ws.Range("I98").Value = 917.6667
$ws.Range("J98").Value = 895
$ws.Range("K98").Value = 2753.0001
$ws.Range("L98").Value = 2685
$ws.Range("M98").Value = -1255.0001
$ws.Range("N98").Value = -5681

# Row 117
$ws.Range("H117").Value = 2245.9
$ws.Range("I117").Value = 957
$ws.Range("J117").Value = 3105.1667
$ws.Range("K117").Value = 2871
$ws.Range("L117").Value = 9315.500100000001
$ws.Range("M117").Value = 571
$ws.Range("N117").Value = -16199.5001

# Row 129
$ws.Range("H129").Value = 2034.1765
$ws.Range("I129").Value = 968.7
$ws.Range("J129").Value = 3556.2856
$ws.Range("K129").Value = 2906.1
$ws.Range("L129").Value = 10668.8568
$ws.Range("M129").Value = 2093.9
$ws.Range("N129").Value = -20668.8568

# Row 137
$ws.Range("H137").Value = 7113.5454
$ws.Range("I137").Value = 4583
$ws.Range("J137").Value = 8062.5
$ws.Range("K137").Value = 13749
$ws.Range("L137").Value = 24187.5
$ws.Range("M137").Value = -8649
$ws.Range("N137").Value = -34387.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2486.3901
$ws.Range("I40").Value = 2363.9062
$ws.Range("J40").Value = 2921.889
$ws.Range("K40").Value = 2363.9062
$ws.Range("L40").Value = 2921.889
$ws.Range("M40").Value = -2227.9062
$ws.Range("N40").Value = -3193.889

# Row 55
$ws.Range("H55").Value = 168.33333
$ws.Range("I55").Value = 152.66667
$ws.Range("J55").Value = 199.66667
$ws.Range("K55").Value = 152.66667
$ws.Range("L55").Value = 199.66667
$ws.Range("M55").Value = 20.33332999999999
$ws.Range("N55").Value = -545.6666700000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 46343.76
$ws.Range("I81").Value = 103776.5
$ws.Range("J81").Value = 8055.2666
$ws.Range("K81").Value = 207553
$ws.Range("L81").Value = 16110.5332
$ws.Range("M81").Value = -206492
$ws.Range("N81").Value = -18232.5332

# Row 84
$ws.Range("H84").Value = 46343.76
$ws.Range("I84").Value = 103776.5
$ws.Range("J84").Value = 8055.2666
$ws.Range("K84").Value = 1037765
$ws.Range("L84").Value = 80552.666
$ws.Range("M84").Value = -1032461
$ws.Range("N84").Value = -91160.666

# Row 109
$ws.Range("H109").Value = 82332.664
$ws.Range("J109").Value = 82332.664
$ws.Range("L109").Value = 82332.664
$ws.Range("N109").Value = -85106.664

# Row 136
$ws.Range("H136").Value = 8658.611999999999
$ws.Range("I136").Value = 9448.105
$ws.Range("J136").Value = 5931.273
$ws.Range("K136").Value = 28344.315
$ws.Range("L136").Value = 17793.819
$ws.Range("M136").Value = -25794.315
$ws.Range("N136").Value = -22893.819

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
